$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '61.409.96'
$ws.Range('E2').Value = '  -5.50%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.093.38'
$ws.Range('E3').Value = '  -8.23%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '506.27'
$ws.Range('E5').Value = '  -3.72%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '166.25'
$ws.Range('E6').Value = '  -9.78%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.578'
$ws.Range('E7').Value = '  -3.35%  '
$ws.Range('E8').Value = '  +0.07%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '3.095.60'
$ws.Range('E9').Value = '  -8.29%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.577'
$ws.Range('E10').Value = '  -7.06%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '50.93'
$ws.Range('E11').Value = '  -10.91%  '
$ws.Range('E12').Value = '  -5.60%  '
$ws.Range('E13').Value = '  -4.92%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '8.66'
$ws.Range('E14').Value = '  -5.98%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.586.53'
$ws.Range('E15').Value = '  -8.11%  '
$ws.Range('E16').Value = '  -8.23%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.092.31'
$ws.Range('E17').Value = '  -8.06%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '61.188.64'
$ws.Range('E18').Value = '  -5.37%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '16.50'
$ws.Range('E19').Value = '  -4.64%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '10.50'
$ws.Range('E20').Value = '  -4.80%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.927'
$ws.Range('E21').Value = '  -3.64%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '353.17'
$ws.Range('E22').Value = '  -4.66%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '78.04'
$ws.Range('E23').Value = '  -3.33%  '
$ws.Range('E24').Value = '  -2.88%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '10.75'
$ws.Range('E25').Value = '  -0.47%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '6.10'
$ws.Range('E26').Value = '  +4.90%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '3.80'
$ws.Range('E27').Value = '  +1.15%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.52'
$ws.Range('E28').Value = '  -4.25%  '
$ws.Range('E29').Value = '  -5.90%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.80'
$ws.Range('E30').Value = '  -7.87%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '623.97'
$ws.Range('E31').Value = '  -5.71%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '27.13'
$ws.Range('E32').Value = '  -7.66%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.15'
$ws.Range('E33').Value = '  -7.88%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '10.89'
$ws.Range('E34').Value = '  -2.12%  '
$ws.Range('E35').Value = '  +0.00%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0999'
$ws.Range('E36').Value = '  -3.85%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '55.65'
$ws.Range('E37').Value = '  -8.84%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '34.81'
$ws.Range('E38').Value = '  -4.34%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.362'
$ws.Range('E39').Value = '  -4.29%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.999'
$ws.Range('E40').Value = '  +0.14%  '
$ws.Range('E41').Value = '  +4.38%  '
$ws.Range('E42').Value = '  -7.18%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.768.56'
$ws.Range('E43').Value = '  -1.95%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.45'
$ws.Range('E44').Value = '  +5.56%  '
$ws.Range('B45').Value = 'WEMIXToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.59'
$ws.Range('E45').Value = '  -0.70%  '
$ws.Range('B46').Value = 'Stacks'
$ws.Range('C46').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.82'
$ws.Range('E46').Value = '  +9.87%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0374'
$ws.Range('E47').Value = '  -3.93%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.87'
$ws.Range('E48').Value = '  +1.42%  '
$ws.Range('E49').Value = '  -10.89%  '
$ws.Range('E50').Value = '  -3.86%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '128.88'
$ws.Range('E51').Value = '  -6.26%  '
